# Mark column N ("in main fun, check that intern fun have the good args")
# with "x" for all_args_here.R and all intern_*.R functions (BACKBONE v10.7 pass),
# matching the rows that already have an "x" in column M.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18)

foreach ($r in $rows) {
    $cell = $ws.Range("N$r")
    $cell.Value = "x"
    $cell.HorizontalAlignment = -4108  # xlCenter
    $cell.VerticalAlignment = -4108    # xlCenter
    $cell.WrapText = $true
}

# Update the active selection shown in the sheet view
$ws.Range("N19").Select()
